$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78:175 down to 79:176
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new record's data.
# (Same as the former row 78 except the date (D) and volume (J) values,
# which reflect the new weekly observation.)
$ws.Range("A78").Value = 3
$ws.Range("B78").Value = "Femacal de La Calera"
$ws.Range("C78").Value = "Coquimbo"
$ws.Range("D78").Value = 44671
$ws.Range("E78").Value = 5
$ws.Range("F78").Value = 100112052
$ws.Range("G78").Value = "Albahaca"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 65
$ws.Range("K78").Value = 4000
$ws.Range("L78").Value = 4000
$ws.Range("M78").Value = 4000
$ws.Range("N78").Value = "$/docena de matas"
$ws.Range("O78").Value = "Provincia de Quillota"
$ws.Range("P78").Value = 667
$ws.Range("Q78").Value = 6
$ws.Range("R78").Value = "Hortaliza"
